$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed cryptos list

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.213.13'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.90%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.857.97'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.79%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7141'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.62%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.47'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.51%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07749'

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3072'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.06%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.10'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.84%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08252'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.22%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.869.70'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.37%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.244'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7166'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.73%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.22'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.22%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.237.20'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.96%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.863'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '244.20'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.51%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007794'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.78%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.15'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.13%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.105.58'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.27%  '

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.04%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.991'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +3.09%  '

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.07%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1595'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.31%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.35'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.34%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.910'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.88%  '

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.25%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.496'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.10%  '

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -3.08%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.397'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.48%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.199'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +2.89%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05185'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.24%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.909'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.11%  '

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.22%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7264'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.35%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.678'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.03%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01855'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.69%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.687'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.15%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.152.36'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.00%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9048'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.148'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.64%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '72.29'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.02%  '

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '101.69'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.60%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.001.35'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.34%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5223'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.64%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.766'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.13%  '

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.29%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.318'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.87%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.864'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.25%  '
